$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '23.193.80'
$c.Style = "Normal"

$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.60%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '1.601.75'
$c.Style = "Normal"

$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.12%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"

$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.05%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '304.85'
$c.Style = "Normal"

$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.84%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '0.3759'
$c.Style = "Normal"

$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.54%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '52.80'
$c.Style = "Normal"

$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = '  +3.94%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.3603'
$c.Style = "Normal"

$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.04%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '1.258'
$c.Style = "Normal"

$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.63%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"

$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.11%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '0.08126'
$c.Style = "Normal"

$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.21%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '22.77'
$c.Style = "Normal"

$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.91%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '6.579'
$c.Style = "Normal"

$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.11%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '7.339'
$c.Style = "Normal"

$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.40%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '0.00001241'
$c.Style = "Normal"

$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.07%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '1.602.91'
$c.Style = "Normal"

$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '93.94'
$c.Style = "Normal"

$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.92%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '0.06921'
$c.Style = "Normal"

$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.07%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '18.10'
$c.Style = "Normal"

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '6.521'
$c.Style = "Normal"

$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.23%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.Style = "Normal"

$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.20%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '12.85'
$c.Style = "Normal"

$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.33%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '23.187.98'
$c.Style = "Normal"

$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.55%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(25, 2)
$c.NumberFormat = "@"
$c.Value = 'Toncoin'
$c.Style = "Normal"

$c = $ws.Cells.Item(25, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c.Style = "Normal"

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '2.424'
$c.Style = "Normal"

$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = '  +2.88%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(26, 2)
$c.NumberFormat = "@"
$c.Value = 'LidoDAOToken'
$c.Style = "Normal"

$c = $ws.Cells.Item(26, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c.Style = "Normal"

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '3.049'
$c.Style = "Normal"

$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = '  +10.39%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '21.13'
$c.Style = "Normal"

$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.21%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '150.70'
$c.Style = "Normal"

$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.02%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '5.267'
$c.Style = "Normal"

$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '134.83'
$c.Style = "Normal"

$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.01%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '2.408'
$c.Style = "Normal"

$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = "@"
$c.Value = '  +2.18%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '6.721'
$c.Style = "Normal"

$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.88%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '1.780.17'
$c.Style = "Normal"

$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.31%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '0.9482'
$c.Style = "Normal"

$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.04%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '0.02763'
$c.Style = "Normal"

$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = "@"
$c.Value = '  +2.15%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '0.07406'
$c.Style = "Normal"

$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.57%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '10.24'
$c.Style = "Normal"

$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.93%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '0.2509'
$c.Style = "Normal"

$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.30%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '6.085'
$c.Style = "Normal"

$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.67%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '0.08740'
$c.Style = "Normal"

$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.88%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '1.400'
$c.Style = "Normal"

$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = '  +3.14%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '0.7074'
$c.Style = "Normal"

$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.70%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '12.38'
$c.Style = "Normal"

$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.00%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '15.85'
$c.Style = "Normal"

$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = '  +4.28%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '0.6503'
$c.Style = "Normal"

$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.11%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '2.319'
$c.Style = "Normal"

$c = $ws.Cells.Item(47, 2)
$c.NumberFormat = "@"
$c.Value = 'PancakeSwap'
$c.Style = "Normal"

$c = $ws.Cells.Item(47, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c.Style = "Normal"

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '4.008'
$c.Style = "Normal"

$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.23%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(48, 2)
$c.NumberFormat = "@"
$c.Value = 'Quant'
$c.Style = "Normal"

$c = $ws.Cells.Item(48, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c.Style = "Normal"

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '133.81'
$c.Style = "Normal"

$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.65%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(49, 2)
$c.NumberFormat = "@"
$c.Value = 'Cronos'
$c.Style = "Normal"

$c = $ws.Cells.Item(49, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c.Style = "Normal"

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '0.07948'
$c.Style = "Normal"

$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.24%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(50, 2)
$c.NumberFormat = "@"
$c.Value = 'Flow'
$c.Style = "Normal"

$c = $ws.Cells.Item(50, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$c.Style = "Normal"

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '1.192'
$c.Style = "Normal"

$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.24%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(51, 2)
$c.NumberFormat = "@"
$c.Value = 'ThetaToken'
$c.Style = "Normal"

$c = $ws.Cells.Item(51, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c.Style = "Normal"

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '1.188'
$c.Style = "Normal"

$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.13%  '
$c.Style = "Normal"
